$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 / 21: "Send PATCH request for lookup value JSON object" -> "Send PATCH request for lookup value"
$ws.Range("H15").Value = "Send PATCH request for lookup value"
$ws.Range("H21").Value = "Send PATCH request for lookup value"

# Row 31: "Send POST request for lookup value JSON object" -> "Send POST request for lookup value"
$ws.Range("H31").Value = "Send POST request for lookup value"

# Row 26: "Send POST request for lookup value JSON object" -> "Send POST request for lookup value " (trailing space retained)
$ws.Range("H26").Value = "Send POST request for lookup value "
$ws.Rows(26).EntireRow.AutoFit()

# Update the active selection to match the author's final cursor position
$ws.Range("I26").Select()
